$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2435.318
$ws.Range("J40").Value = 2463.8
$ws.Range("L40").Value = 2463.8
$ws.Range("N40").Value = -2813.8

$ws.Range("H64").Value = 37711
$ws.Range("J64").Value = 5369
$ws.Range("L64").Value = 5369
$ws.Range("N64").Value = -5865

$ws.Range("H67").Value = 37711
$ws.Range("J67").Value = 5369
$ws.Range("L67").Value = 5369
$ws.Range("N67").Value = -7085

$ws.Range("H70").Value = 7410438
$ws.Range("I70").Value = 22226154
$ws.Range("J70").Value = 2579.5557
$ws.Range("K70").Value = 66678462
$ws.Range("L70").Value = 7738.6671
$ws.Range("M70").Value = -66678192
$ws.Range("N70").Value = -8278.667099999999

$ws.Range("H73").Value = 7410438
$ws.Range("I73").Value = 22226154
$ws.Range("J73").Value = 2579.5557
$ws.Range("K73").Value = 66678462
$ws.Range("L73").Value = 7738.6671
$ws.Range("M73").Value = -66677526
$ws.Range("N73").Value = -9610.667099999999

$ws.Range("H94").Value = 55698810
$ws.Range("I94").Value = 83343220
$ws.Range("K94").Value = 83343220
$ws.Range("M94").Value = -83342769

$ws.Range("H129").Value = 66668428
$ws.Range("I129").Value = 1437.8889
$ws.Range("K129").Value = 4313.6667
$ws.Range("M129").Value = 686.3333000000002

$ws.Range("H135").Value = 3490.125
$ws.Range("I135").Value = 4254.909
$ws.Range("K135").Value = 38294.181
$ws.Range("M135").Value = -35759.181

$ws.Range("H138").Value = 8986.322
$ws.Range("I138").Value = 8402.25
$ws.Range("J138").Value = 9355.210999999999
$ws.Range("K138").Value = 25206.75
$ws.Range("L138").Value = 28065.633
$ws.Range("M138").Value = -20066.75
$ws.Range("N138").Value = -38345.633

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2966.6223
$ws.Range("I32").Value = 2893.0476
$ws.Range("K32").Value = 2893.0476
$ws.Range("M32").Value = -2606.0476

$ws.Range("H43").Value = 14570.8
$ws.Range("I43").Value = 19868.5
$ws.Range("J43").Value = 13246.375
$ws.Range("K43").Value = 19868.5
$ws.Range("L43").Value = 13246.375
$ws.Range("M43").Value = -19555.5
$ws.Range("N43").Value = -13872.375

$ws.Range("H61").Value = 11061.883
$ws.Range("I61").Value = 13087.917
$ws.Range("K61").Value = 13087.917
$ws.Range("M61").Value = -12875.917

$ws.Range("H136").Value = 11061.883
$ws.Range("I136").Value = 13087.917
$ws.Range("K136").Value = 39263.751
$ws.Range("M136").Value = -36713.751

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 8871.632
$ws.Range("I94").Value = 10862.143
$ws.Range("K94").Value = 10862.143
$ws.Range("M94").Value = -10411.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2414.6924
$ws.Range("I16").Value = 3317
$ws.Range("J16").Value = 1641.2858
$ws.Range("K16").Value = 3317
$ws.Range("L16").Value = 1641.2858
$ws.Range("M16").Value = -3030
$ws.Range("N16").Value = -2215.2858

$ws.Range("H58").Value = 7402.4287
$ws.Range("J58").Value = 3874.6155
$ws.Range("L58").Value = 3874.6155
$ws.Range("N58").Value = -4280.6155

$ws.Range("H113").Value = 2414.6924
$ws.Range("I113").Value = 3317
$ws.Range("J113").Value = 1641.2858
$ws.Range("K113").Value = 3317
$ws.Range("L113").Value = 1641.2858
$ws.Range("M113").Value = -1147
$ws.Range("N113").Value = -5981.2858

$ws.Range("H122").Value = 2012
$ws.Range("I122").Value = 2012
$ws.Range("K122").Value = 6036
$ws.Range("M122").Value = -3586

$ws.Range("H134").Value = 1791295.5
$ws.Range("I134").Value = 2506386.2
$ws.Range("K134").Value = 7519158.600000001
$ws.Range("M134").Value = -7516623.600000001

$ws.Range("H136").Value = 7402.4287
$ws.Range("J136").Value = 3874.6155
$ws.Range("L136").Value = 11623.8465
$ws.Range("N136").Value = -16723.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 280.7619
$ws.Range("I6").Value = 301.64285
$ws.Range("J6").Value = 239
$ws.Range("K6").Value = 904.9285500000001
$ws.Range("L6").Value = 717
$ws.Range("M6").Value = -791.9285500000001
$ws.Range("N6").Value = -943

$ws.Range("H8").Value = 383.0909
$ws.Range("I8").Value = 383.0909
$ws.Range("K8").Value = 1149.2727
$ws.Range("M8").Value = -1010.2727

$ws.Range("H33").Value = 266
$ws.Range("I33").Value = 166.83333
$ws.Range("J33").Value = 414.75
$ws.Range("K33").Value = 1000.99998
$ws.Range("L33").Value = 2488.5
$ws.Range("M33").Value = -717.9999799999999
$ws.Range("N33").Value = -3054.5

$ws.Range("H44").Value = 1044.1111
$ws.Range("I44").Value = 527.38464
$ws.Range("K44").Value = 1582.15392
$ws.Range("M44").Value = -1184.15392

$ws.Range("H51").Value = 557562.9
$ws.Range("I51").Value = 834141.9399999999
$ws.Range("K51").Value = 2502425.82
$ws.Range("M51").Value = -2501965.82

$ws.Range("H63").Value = 2392.3333
$ws.Range("I63").Value = 2091.5
$ws.Range("J63").Value = 2994
$ws.Range("K63").Value = 6274.5
$ws.Range("L63").Value = 8982
$ws.Range("M63").Value = -5525.5
$ws.Range("N63").Value = -10480

$ws.Range("H66").Value = 2392.3333
$ws.Range("I66").Value = 2091.5
$ws.Range("J66").Value = 2994
$ws.Range("K66").Value = 18823.5
$ws.Range("L66").Value = 26946
$ws.Range("M66").Value = -15079.5
$ws.Range("N66").Value = -34434

$ws.Range("H103").Value = 2990.25
$ws.Range("I103").Value = 2498.8
$ws.Range("J103").Value = 3809.3333
$ws.Range("K103").Value = 7496.400000000001
$ws.Range("L103").Value = 11427.9999
$ws.Range("M103").Value = -6617.400000000001
$ws.Range("N103").Value = -13185.9999

$ws.Range("H107").Value = 1343.4166
$ws.Range("I107").Value = 575
$ws.Range("J107").Value = 1497.1
$ws.Range("K107").Value = 1725
$ws.Range("L107").Value = 4491.299999999999
$ws.Range("M107").Value = 195
$ws.Range("N107").Value = -8331.299999999999

$ws.Range("H114").Value = 10927.7
$ws.Range("I114").Value = 768
$ws.Range("J114").Value = 34633.668
$ws.Range("K114").Value = 2304
$ws.Range("L114").Value = 103901.004
$ws.Range("M114").Value = 950
$ws.Range("N114").Value = -110409.004

$ws.Range("H117").Value = 8132.8335
$ws.Range("I117").Value = 2695.4285
$ws.Range("J117").Value = 11593
$ws.Range("K117").Value = 8086.2855
$ws.Range("L117").Value = 34779
$ws.Range("M117").Value = -4644.2855
$ws.Range("N117").Value = -41663

$ws.Range("H129").Value = 20834852
$ws.Range("I129").Value = 1428.2858
$ws.Range("J129").Value = 37038628
$ws.Range("K129").Value = 4284.857400000001
$ws.Range("L129").Value = 111115884
$ws.Range("M129").Value = 715.1425999999992
$ws.Range("N129").Value = -111125884

$ws.Range("H137").Value = 7461.1904
$ws.Range("J137").Value = 12394.5
$ws.Range("L137").Value = 37183.5
$ws.Range("N137").Value = -47383.5

$ws.Range("H139").Value = 2310892.8
$ws.Range("I139").Value = 4287258.5
$ws.Range("K139").Value = 12861775.5
$ws.Range("M139").Value = -12856635.5

$ws.Range("H140").Value = 1727.909
$ws.Range("I140").Value = 1475.8
$ws.Range("K140").Value = 4427.4
$ws.Range("M140").Value = 752.6000000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4447.8
$ws.Range("J80").Value = 3335.9333
$ws.Range("L80").Value = 3335.9333
$ws.Range("N80").Value = -5331.933300000001

$ws.Range("H83").Value = 4447.8
$ws.Range("J83").Value = 3335.9333
$ws.Range("L83").Value = 16679.6665
$ws.Range("N83").Value = -26663.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3796.373
$ws.Range("I22").Value = 5066.4414
$ws.Range("J22").Value = 2487.818
$ws.Range("K22").Value = 5066.4414
$ws.Range("L22").Value = 2487.818
$ws.Range("M22").Value = -4771.4414
$ws.Range("N22").Value = -3077.818

$ws.Range("H27").Value = 3796.373
$ws.Range("I27").Value = 5066.4414
$ws.Range("J27").Value = 2487.818
$ws.Range("K27").Value = 5066.4414
$ws.Range("L27").Value = 2487.818
$ws.Range("M27").Value = -4959.4414
$ws.Range("N27").Value = -2701.818

$ws.Range("H55").Value = 564.63635
$ws.Range("I55").Value = 563.12
$ws.Range("K55").Value = 563.12
$ws.Range("M55").Value = -390.12

$ws.Range("H139").Value = 127749.5
$ws.Range("J139").Value = 98666
$ws.Range("L139").Value = 98666
$ws.Range("N139").Value = -108946

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9041.083000000001
$ws.Range("I122").Value = 7928.2856
$ws.Range("J122").Value = 10599
$ws.Range("K122").Value = 23784.8568
$ws.Range("L122").Value = 31797
$ws.Range("M122").Value = -21334.8568

$ws.Range("H126").Value = 28880.588
$ws.Range("I126").Value = 69515.664
$ws.Range("J126").Value = 6716
$ws.Range("K126").Value = 208546.992
$ws.Range("L126").Value = 20148
$ws.Range("M126").Value = -206076.992
$ws.Range("N126").Value = -25088
